$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: ORD_MGR_02 - shorten title, swap the "steps" / "sample data" columns
$ws.Cells.Item(2, 2).Value = "Lọc đơn hàng"
$ws.Cells.Item(2, 3).Value = "Status='Pending'"
$ws.Cells.Item(2, 4).Value = "1. Action='List'`n2. Status='Pending'"

# Row 3: ORD_MGR_01 - swap the "steps" / "sample data" columns
$ws.Cells.Item(3, 3).Value = "Status=null"
$ws.Cells.Item(3, 4).Value = "1. Action='List'`n2. Status=null"

# Row 4: was ORD_MGR_04 / Xoa don hang -> now ORD_MGR_06 (same title), swap data/steps
$ws.Cells.Item(4, 1).Value = "ORD_MGR_06"
$ws.Cells.Item(4, 2).Value = "Xóa đơn hàng"
$ws.Cells.Item(4, 3).Value = "ID=5"
$ws.Cells.Item(4, 4).Value = "1. Action='Delete'`n2. ID=5"
$ws.Cells.Item(4, 5).Value = "Call delete -> Redirect"
$ws.Cells.Item(4, 6).Value = "OK"
$ws.Cells.Item(4, 7).Value = "PASS"

# Row 5: was ORD_MGR_03 / Cap nhat trang thai don -> now ORD_MGR_04 / Luu cap nhat don hang
$ws.Cells.Item(5, 1).Value = "ORD_MGR_04"
$ws.Cells.Item(5, 2).Value = "Lưu cập nhật đơn hàng"
$ws.Cells.Item(5, 3).Value = "ID=10, Total=500k"
$ws.Cells.Item(5, 4).Value = "1. Action='SaveOrUpdate'`n2. Params đầy đủ"
$ws.Cells.Item(5, 5).Value = "Call Service Save -> Redirect List"
$ws.Cells.Item(5, 6).Value = "OK"
$ws.Cells.Item(5, 7).Value = "PASS"

# Row 6 (new): ORD_MGR_08 / Loi he thong
$ws.Cells.Item(6, 1).Value = "ORD_MGR_08"
$ws.Cells.Item(6, 2).Value = "Lỗi hệ thống"
$ws.Cells.Item(6, 3).Value = "Exception"
$ws.Cells.Item(6, 4).Value = "Service ném lỗi"
$ws.Cells.Item(6, 5).Value = "Forward trang error.jsp"
$ws.Cells.Item(6, 6).Value = "OK"
$ws.Cells.Item(6, 7).Value = "PASS"

# Row 7 (new): ORD_MGR_05 / Cap nhat trang thai nhanh
$ws.Cells.Item(7, 1).Value = "ORD_MGR_05"
$ws.Cells.Item(7, 2).Value = "Cập nhật trạng thái nhanh"
$ws.Cells.Item(7, 3).Value = "ID=10, St='Done'"
$ws.Cells.Item(7, 4).Value = "1. Action='UpdateStatus'`n2. ID=10, Status='Done'"
$ws.Cells.Item(7, 5).Value = "Call update -> Redirect"
$ws.Cells.Item(7, 6).Value = "OK"
$ws.Cells.Item(7, 7).Value = "PASS"

# Row 8 (new): ORD_MGR_03 / Hien form sua don
$ws.Cells.Item(8, 1).Value = "ORD_MGR_03"
$ws.Cells.Item(8, 2).Value = "Hiện form sửa đơn"
$ws.Cells.Item(8, 3).Value = "ID=5"
$ws.Cells.Item(8, 4).Value = "1. Action='AddOrEdit', ID=5"
$ws.Cells.Item(8, 5).Value = "Forward OrdersManager.jsp"
$ws.Cells.Item(8, 6).Value = "OK"
$ws.Cells.Item(8, 7).Value = "PASS"

# Row 9 (new): ORD_MGR_07 / Action Null -> List
$ws.Cells.Item(9, 1).Value = "ORD_MGR_07"
$ws.Cells.Item(9, 2).Value = "Action Null -> List"
$ws.Cells.Item(9, 3).Value = "Null"
$ws.Cells.Item(9, 4).Value = "Action=null"
$ws.Cells.Item(9, 5).Value = "Mặc định gọi getAllOrders"
$ws.Cells.Item(9, 6).Value = "OK"
$ws.Cells.Item(9, 7).Value = "PASS"

# Give new rows 6-9 the same "PASS" (green bold) formatting already used in column G
$ws.Range("G5").Copy()
$ws.Range("G6:G9").PasteSpecial(-4122)

# Resize columns B:E to match the recalculated best-fit widths for the new content
$ws.Columns.Item(2).ColumnWidth = 23.5
$ws.Columns.Item(3).ColumnWidth = 16.5
$ws.Columns.Item(4).ColumnWidth = 24.666666666666668
$ws.Columns.Item(5).ColumnWidth = 28.833333333333332
